# Apply the "add text to pdf screenshot" edit to the "comment" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comment")

# Set D2:D4 to verbosity level 3
$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 3
$ws.Range("D4").Value = 3

# Make this sheet active and move the selection to D4
$ws.Activate()
$ws.Range("D4").Select()
